$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 96 and 97, shifting the existing data
# (previously rows 96-128) down to rows 98-130.
$ws.Rows("96:97").Insert()

# Populate the two newly inserted rows with the new "Albaricoque" records.
# Columns A,B,C,E,F,G,H,I,J carry the same constant values used throughout
# the whole data range.

$ws.Range("A96").Value = 9
$ws.Range("B96").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C96").Value = "Metropolitana"
$ws.Range("D96").Value = 44907
$ws.Range("E96").Value = 13
$ws.Range("F96").Value = "Fruta"
$ws.Range("G96").Value = 100103
$ws.Range("H96").Value = "Frutos de hueso (carozo)"
$ws.Range("I96").Value = 100103003
$ws.Range("J96").Value = "Damasco"
$ws.Range("K96").Value = "Albaricoque"
$ws.Range("L96").Value = "Primera"
$ws.Range("M96").Value = 280
$ws.Range("N96").Value = 13000
$ws.Range("O96").Value = 13000
$ws.Range("P96").Value = 13000
$ws.Range("Q96").Value = "$/caja 16 kilos granel"
$ws.Range("R96").Value = "Región de O'Higgins"
$ws.Range("S96").Value = 812
$ws.Range("T96").Value = 16

$ws.Range("A97").Value = 9
$ws.Range("B97").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C97").Value = "Metropolitana"
$ws.Range("D97").Value = 44907
$ws.Range("E97").Value = 13
$ws.Range("F97").Value = "Fruta"
$ws.Range("G97").Value = 100103
$ws.Range("H97").Value = "Frutos de hueso (carozo)"
$ws.Range("I97").Value = 100103003
$ws.Range("J97").Value = "Damasco"
$ws.Range("K97").Value = "Albaricoque"
$ws.Range("L97").Value = "Segunda"
$ws.Range("M97").Value = 300
$ws.Range("N97").Value = 11000
$ws.Range("O97").Value = 11000
$ws.Range("P97").Value = 11000
$ws.Range("Q97").Value = "$/caja 16 kilos granel"
$ws.Range("R97").Value = "Región de O'Higgins"
$ws.Range("S97").Value = 688
$ws.Range("T97").Value = 16
